$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "7"

# Create sheet "8" as an exact copy of sheet "7" while it still only has the
# original header + row2, so the new sheet inherits the same column widths
# (xl/cols) and styles without us having to reconstruct them by hand.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "8"

# ---------------------------------------------------------------------
# Sheet "7": extend row 2 and add rows 3-4 (3 data rows total), then
# merge A2:A4 and M2:M3, matching the target sheet1.xml
# ---------------------------------------------------------------------

$ws1.Range("A2:A4").Merge()
$ws1.Range("M2:M3").Merge()

$ws1.Range("A2").Value = "4"
$ws1.Range("B2").Value = "VNDUSD"
$ws1.Range("C2").Value = "14:57:30.000793"
$ws1.Range("D2").Value = "Sell"
$ws1.Range("E2").Value = "link 1"
$ws1.Range("F2").Value = "link 2"
$ws1.Range("G2").Value = "link 3"
$ws1.Range("H2").Value = "link 4"
$ws1.Range("I2").Value = "link 5"
$ws1.Range("J2").Value = 3
$ws1.Range("K2").Value = "this is my comment"
$ws1.Range("L2").Value = 601
$ws1.Range("M2").Formula = "=SUM(J2:J4)"

$ws1.Range("B3").Value = "PAIR2"
$ws1.Range("C3").Value = "14:57:56.200590"
$ws1.Range("D3").Value = "Sell"
$ws1.Range("E3").Value = "link 1"
$ws1.Range("F3").Value = "link 2"
$ws1.Range("G3").Value = "link 3"
$ws1.Range("H3").Value = "link 4"
$ws1.Range("I3").Value = "link 5"
$ws1.Range("J3").Value = 3
$ws1.Range("K3").Value = "this is my comment"
$ws1.Range("L3").Value = 81

$ws1.Range("B4").Value = "PAIR2"
$ws1.Range("C4").Value = "14:58:01.136710"
$ws1.Range("D4").Value = "Sell"
$ws1.Range("E4").Value = "link 1"
$ws1.Range("F4").Value = "link 2"
$ws1.Range("G4").Value = "link 3"
$ws1.Range("H4").Value = "link 4"
$ws1.Range("I4").Value = "link 5"
$ws1.Range("J4").Value = 3
$ws1.Range("K4").Value = "this is my comment"
$ws1.Range("L4").Value = 801

# Re-apply the border + centered alignment (style index 1 in the original
# file) across the newly added rows so every cell - including the blank
# ones left behind by the merges - carries the same cell style as the rest
# of the table.
$rng1 = $ws1.Range("A2:M4")
$rng1.Borders.LineStyle = 1
$rng1.HorizontalAlignment = -4108
$rng1.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Sheet "8": overwrite the copied row 2 and add row 3 (2 data rows total),
# then merge A2:A3 and M2:M3, matching the target sheet2.xml
# ---------------------------------------------------------------------

$ws2.Range("A2:A3").Merge()
$ws2.Range("M2:M3").Merge()

$ws2.Range("A2").Value = "4"
$ws2.Range("B2").Value = "VNDUSD"
$ws2.Range("C2").Value = "14:56:57.405001"
$ws2.Range("D2").Value = "Sell"
$ws2.Range("E2").Value = "link 1"
$ws2.Range("F2").Value = "link 2"
$ws2.Range("G2").Value = "link 3"
$ws2.Range("H2").Value = "link 4"
$ws2.Range("I2").Value = "link 5"
$ws2.Range("J2").Value = 3
$ws2.Range("K2").Value = "this is my comment"
$ws2.Range("L2").Value = 201
$ws2.Range("M2").Formula = "=SUM(J2:J3)"

$ws2.Range("B3").Value = "VNDUSD"
$ws2.Range("C3").Value = "14:57:03.901108"
$ws2.Range("D3").Value = "Sell"
$ws2.Range("E3").Value = "link 1"
$ws2.Range("F3").Value = "link 2"
$ws2.Range("G3").Value = "link 3"
$ws2.Range("H3").Value = "link 4"
$ws2.Range("I3").Value = "link 5"
$ws2.Range("J3").Value = 3
$ws2.Range("K3").Value = "this is my comment"
$ws2.Range("L3").Value = 401

$rng2 = $ws2.Range("A2:M3")
$rng2.Borders.LineStyle = 1
$rng2.HorizontalAlignment = -4108
$rng2.VerticalAlignment = -4108

# Restore sheet "7" as the selected/active tab, matching the original file
$ws1.Activate()
